$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.074.69"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.993.25"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.68"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.22"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.06%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.988.61"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.75"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.33%  "
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.496.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.02"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.100.99"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.999.29"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "455.22"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.90"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.685"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.58"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -9.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.14"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.34%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.72%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.63"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.91"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.57"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0800"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.71"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.09"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.13"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.14"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.122"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.86"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -11.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "390.40"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0355"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.266"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -7.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.725.44"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.31"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.41"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.108"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.17"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.39%  "
